# Generate Report for Handback
#
# This script updates the localization-status workbook to reflect a
# completed handback: the "Status" column moves from "Ready for handoff"
# to "Handed back: in sync with en-US", the per-language sheets get their
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns populated (with a hyperlink on the target-file cell), and a
# handful of columns are widened so the new content is not clipped.

$wb  = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shown on Overview!E2, Overview!F2, zh-cn!C2 and de-de!C2)
# ---------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"
$ovw.Range("E2").Value2 = $newStatus
$ovw.Range("F2").Value2 = $newStatus
$zh.Range("C2").Value2  = $newStatus
$de.Range("C2").Value2  = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File (H2, hyperlinked like A2),
#    Latest Handback File (I2) and Latest Handback DateTime (J2)
# ---------------------------------------------------------------------
$mdName  = "1c6c813b-57cb-4dd7-8583-35a97b77e7ce.md"
$mdUrl   = "https://github.com/OpenLocalizationTestOrg/oltest/blob/a0f8b315d86811746969df6e388f00ff75af2dbb/e2e/1c6c813b-57cb-4dd7-8583-35a97b77e7ce.md"

$zh.Hyperlinks.Add($zh.Range("H2"), $mdUrl, "", "", $mdName) | Out-Null
$zh.Range("I2").Value2 = "1c6c813b-57cb-4dd7-8583-35a97b77e7ce.038f658be19e0eb48ad776c09a83d71ac2a0fb15.zh-cn.xlf"
$zh.Range("J2").Value2 = "2016-07-04 09:10:28"

# ---------------------------------------------------------------------
# 3. de-de sheet: same treatment, but with its own handback datetime
# ---------------------------------------------------------------------
$de.Hyperlinks.Add($de.Range("H2"), $mdUrl, "", "", $mdName) | Out-Null
$de.Range("I2").Value2 = "1c6c813b-57cb-4dd7-8583-35a97b77e7ce.038f658be19e0eb48ad776c09a83d71ac2a0fb15.de-de.xlf"
$de.Range("J2").Value2 = "2016-07-04 09:10:43"

# ---------------------------------------------------------------------
# 4. Widen columns so the longer status text / hyperlink text fit
# ---------------------------------------------------------------------
# Overview: zh-cn (E) and de-de (F) status columns
$ovw.Columns.Item(5).ColumnWidth = 29.166666666666668
$ovw.Columns.Item(6).ColumnWidth = 29.166666666666668

foreach ($ws in @($zh, $de)) {
    $ws.Columns.Item(3).ColumnWidth  = 29.166666666666668   # C - Status
    $ws.Columns.Item(8).ColumnWidth  = 39.166666666666664   # H - Latest Target File
    $ws.Columns.Item(9).ColumnWidth  = 39.166666666666664   # I - Latest Handback File
}
